$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: the "Além disso ..." sentence used to be split into three
# runs with a spurious <w:proofErr> pair wrapping the lone word "e"
# (an artifact of Word's grammar checker). Re-merge it into one
# contiguous sentence / run, exactly as a user accepting the grammar
# suggestion (or simply retyping the sentence) would in Word.
# ---------------------------------------------------------------------
$old1 = "Além disso o que o Terraform cria faz parte da configuração do Terraform que e o terraform.lock.hcl"
$new1 = "Além disso o que o Terraform cria faz parte da configuração do Terraform que e o terraform.lock.hcl"

$found1 = $d.Content.Find.Execute(
    $old1, $true, $false, $false, $false, $false,
    $true, 1, $false, $new1, 2)

# ---------------------------------------------------------------------
# Change 2: add a new line "-Existem outros comandos " right after the
# "Terraform destroy - Destroy previously-created infrastructure "
# paragraph (and before the blank paragraph that follows it).
# ---------------------------------------------------------------------
$destroyIdx = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i++
    if ($p.Range.Text -like "*Destroy previously-created infrastructure*") {
        $destroyIdx = $i
    }
}

$blankAfter = $d.Paragraphs($destroyIdx + 1)
$blankAfter.Range.InsertBefore("-Existem outros comandos `r")
